$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and the Mantle/Filecoin row swap)
# D-column values that look numeric are forced to Text format first so Excel
# does not silently convert them from strings to numbers.

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '68.341.27'
$ws.Cells.Item(2, 5).Value = '  +0.68%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.846.68'
$ws.Cells.Item(3, 5).Value = '  -0.14%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.04%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '601.65'
$ws.Cells.Item(5, 5).Value = '  +0.54%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '172.09'
$ws.Cells.Item(6, 5).Value = '  +3.42%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.846.40'
$ws.Cells.Item(7, 5).Value = '  -0.10%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.528'
$ws.Cells.Item(9, 5).Value = '  +0.65%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.82%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.54'
$ws.Cells.Item(11, 5).Value = '  +3.32%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +1.38%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +14.99%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '37.07'
$ws.Cells.Item(14, 5).Value = '  +0.77%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.493.54'
$ws.Cells.Item(15, 5).Value = '  -0.14%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.840.36'
$ws.Cells.Item(16, 5).Value = '  -0.80%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '68.410.96'
$ws.Cells.Item(17, 5).Value = '  +0.73%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '18.37'
$ws.Cells.Item(18, 5).Value = '  +1.52%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.45'
$ws.Cells.Item(19, 5).Value = '  +1.43%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.72%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.92'
$ws.Cells.Item(21, 5).Value = '  -0.25%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '470.49'
$ws.Cells.Item(22, 5).Value = '  +1.49%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.733'
$ws.Cells.Item(23, 5).Value = '  +0.66%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -2.65%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '83.72'
$ws.Cells.Item(25, 5).Value = '  +0.60%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +1.97%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.19'
$ws.Cells.Item(27, 5).Value = '  +0.68%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.55'
$ws.Cells.Item(28, 5).Value = '  +5.76%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.25%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.001.42'
$ws.Cells.Item(31, 5).Value = '  -0.04%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.77'
$ws.Cells.Item(32, 5).Value = '  +0.36%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.31'
$ws.Cells.Item(33, 5).Value = '  -0.06%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '31.16'
$ws.Cells.Item(34, 5).Value = '  +0.54%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '9.40'
$ws.Cells.Item(35, 5).Value = '  +0.85%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '3.810.86'
$ws.Cells.Item(36, 5).Value = '  -0.48%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.94'
$ws.Cells.Item(37, 5).Value = '  +20.81%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.105'
$ws.Cells.Item(38, 5).Value = '  +1.18%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Filecoin'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.98'
$ws.Cells.Item(39, 5).Value = '  +1.57%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Mantle'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.02'
$ws.Cells.Item(40, 5).Value = '  +0.44%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.33%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.02%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.320'
$ws.Cells.Item(43, 5).Value = '  +2.59%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.99'
$ws.Cells.Item(44, 5).Value = '  +0.70%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.02%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '8.76'
$ws.Cells.Item(46, 5).Value = '  +3.11%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '418.96'
$ws.Cells.Item(47, 5).Value = '  -1.83%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +8.30%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '46.75'
$ws.Cells.Item(49, 5).Value = '  -0.80%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '142.07'
$ws.Cells.Item(50, 5).Value = '  -0.89%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +1.70%  '
